{"js": "// Update the answer table's 100 arithmetic results (20 rows x 5 cols) to the\n// new set of expressions, preserving each cell's existing formatting.\nconst newValues = [\n  [\"7+85=92\", \"19+62=81\", \"46-17=29\", \"13-8=5\", \"33-8=25\"],\n  [\"54+18=72\", \"66+7=73\", \"8+77=85\", \"87-49=38\", \"9+49=58\"],\n  [\"57-8=49\", \"71-52=19\", \"64+29=93\", \"17+44=61\", \"9+39=48\"],\n  [\"45-39=6\", \"49+26=75\", \"21-13=8\", \"82-38=44\", \"39+55=94\"],\n  [\"93-18=75\", \"3+78=81\", \"74-16=58\", \"41-15=26\", \"64-16=48\"],\n  [\"81-57=24\", \"8+18=26\", \"44+8=52\", \"52-38=14\", \"68+7=75\"],\n  [\"47+35=82\", \"66+17=83\", \"82-8=74\", \"26+9=35\", \"57+4=61\"],\n  [\"46+15=61\", \"56+9=65\", \"18+55=73\", \"95-67=28\", \"8+7=15\"],\n  [\"25+69=94\", \"25+67=92\", \"66+17=83\", \"81-68=13\", \"71-7=64\"],\n  [\"18+43=61\", \"59+32=91\", \"9+34=43\", \"73-34=39\", \"39+48=87\"],\n  [\"72-19=53\", \"19+23=42\", \"68-49=19\", \"26+39=65\", \"57+38=95\"],\n  [\"47-18=29\", \"78-49=29\", \"35+37=72\", \"64+8=72\", \"34-8=26\"],\n  [\"62-4=58\", \"18+64=82\", \"95-49=46\", \"79+4=83\", \"49+43=92\"],\n  [\"69+29=98\", \"84-49=35\", \"37+37=74\", \"7+86=93\", \"69+14=83\"],\n  [\"60-23=37\", \"3+89=92\", \"50-18=32\", \"66-47=19\", \"44-28=16\"],\n  [\"42-38=4\", \"6+66=72\", \"22+49=71\", \"68+19=87\", \"30-14=16\"],\n  [\"83-65=18\", \"12-8=4\", \"4+7=11\", \"5+28=33\", \"93-6=87\"],\n  [\"57+28=85\", \"59+37=96\", \"90-31=59\", \"58+29=87\", \"19+5=24\"],\n  [\"47+15=62\", \"25+67=92\", \"92-18=74\", \"52-14=38\", \"71-53=18\"],\n  [\"7+9=16\", \"19+28=47\", \"29+38=67\", \"23-7=16\", \"38+49=87\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the answer table's 100 arithmetic results (20 rows x 5 cols) to the\n# new set of expressions, preserving each cell's existing formatting.\n$newValues = @(\n    ,@(\"7+85=92\", \"19+62=81\", \"46-17=29\", \"13-8=5\", \"33-8=25\")\n    ,@(\"54+18=72\", \"66+7=73\", \"8+77=85\", \"87-49=38\", \"9+49=58\")\n    ,@(\"57-8=49\", \"71-52=19\", \"64+29=93\", \"17+44=61\", \"9+39=48\")\n    ,@(\"45-39=6\", \"49+26=75\", \"21-13=8\", \"82-38=44\", \"39+55=94\")\n    ,@(\"93-18=75\", \"3+78=81\", \"74-16=58\", \"41-15=26\", \"64-16=48\")\n    ,@(\"81-57=24\", \"8+18=26\", \"44+8=52\", \"52-38=14\", \"68+7=75\")\n    ,@(\"47+35=82\", \"66+17=83\", \"82-8=74\", \"26+9=35\", \"57+4=61\")\n    ,@(\"46+15=61\", \"56+9=65\", \"18+55=73\", \"95-67=28\", \"8+7=15\")\n    ,@(\"25+69=94\", \"25+67=92\", \"66+17=83\", \"81-68=13\", \"71-7=64\")\n    ,@(\"18+43=61\", \"59+32=91\", \"9+34=43\", \"73-34=39\", \"39+48=87\")\n    ,@(\"72-19=53\", \"19+23=42\", \"68-49=19\", \"26+39=65\", \"57+38=95\")\n    ,@(\"47-18=29\", \"78-49=29\", \"35+37=72\", \"64+8=72\", \"34-8=26\")\n    ,@(\"62-4=58\", \"18+64=82\", \"95-49=46\", \"79+4=83\", \"49+43=92\")\n    ,@(\"69+29=98\", \"84-49=35\", \"37+37=74\", \"7+86=93\", \"69+14=83\")\n    ,@(\"60-23=37\", \"3+89=92\", \"50-18=32\", \"66-47=19\", \"44-28=16\")\n    ,@(\"42-38=4\", \"6+66=72\", \"22+49=71\", \"68+19=87\", \"30-14=16\")\n    ,@(\"83-65=18\", \"12-8=4\", \"4+7=11\", \"5+28=33\", \"93-6=87\")\n    ,@(\"57+28=85\", \"59+37=96\", \"90-31=59\", \"58+29=87\", \"19+5=24\")\n    ,@(\"47+15=62\", \"25+67=92\", \"92-18=74\", \"52-14=38\", \"71-53=18\")\n    ,@(\"7+9=16\", \"19+28=47\", \"29+38=67\", \"23-7=16\", \"38+49=87\")\n)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $row = $newValues[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $table.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
